# nbody/Results/Analysis.xlsx - "Updated results and analysis - hard coded 8 in the
# OpenMP statement threw the analysis and results off a bit"
#
# The raw timing results (CPU threads = 2/4/8, for particle counts 512..32768)
# were re-measured after fixing a hard-coded thread count of 8 in the OpenMP
# statement. This updates the raw timing table (C14:E20) and its duplicate
# "Average" row (row 7, pulled in via the delete-named-range import), lets the
# dependent Speedup/Efficiency formulas (N13:P19, S13:U19) recalc naturally,
# switches the three derived line charts to smoothed lines, and restores the
# view (zoom + selection) to how the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Raw CPU timing table (threads=2/4/8 columns C/D/E, rows 14-20 = sizes
#     512,1024,2048,4096,8192,16384,32768) ---------------------------------
$ws.Range("C14").Value = 6
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 4

$ws.Range("C15").Value = 26
$ws.Range("D15").Value = 14
$ws.Range("E15").Value = 10

$ws.Range("C16").Value = 101
$ws.Range("D16").Value = 54
$ws.Range("E16").Value = 40

$ws.Range("C17").Value = 405
$ws.Range("D17").Value = 211
$ws.Range("E17").Value = 140

$ws.Range("C18").Value = 1614
$ws.Range("D18").Value = 835
$ws.Range("E18").Value = 542

$ws.Range("C19").Value = 6469
$ws.Range("D19").Value = 3376
$ws.Range("E19").Value = 2175

$ws.Range("C20").Value = 25910
$ws.Range("D20").Value = 14009
$ws.Range("E20").Value = 8652

# --- Row 7 ("Average" row) mirrors the same raw results under the
#     cpu_data_<size>_<threads> named columns -------------------------------
$ws.Range("D7").Value = 211
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 542
$ws.Range("G7").Value = 26
$ws.Range("H7").Value = 25910
$ws.Range("I7").Value = 3376
$ws.Range("K7").Value = 2175
$ws.Range("L7").Value = 835
$ws.Range("M7").Value = 8652
$ws.Range("N7").Value = 101
$ws.Range("O7").Value = 54
$ws.Range("P7").Value = 6
$ws.Range("S7").Value = 1614
$ws.Range("T7").Value = 3
$ws.Range("U7").Value = 140
$ws.Range("Y7").Value = 10
$ws.Range("Z7").Value = 40
$ws.Range("AA7").Value = 14009
$ws.Range("AB7").Value = 405
$ws.Range("AC7").Value = 6469

# --- Switch the CPU Time / Speedup / Efficiency line charts to smoothed
#     lines (GPU bar chart, the first ChartObject, is untouched) -----------
for ($i = 2; $i -le $ws.ChartObjects().Count; $i++) {
    $chart = $ws.ChartObjects($i).Chart
    for ($j = 1; $j -le $chart.SeriesCollection().Count; $j++) {
        $chart.SeriesCollection($j).Smooth = $true
    }
}

# --- Restore the view: zoom 85% -> 70%, drop the scrolled topLeftCell, and
#     move the selection to O36 -------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$win.Zoom = 70
$ws.Range("O36").Select()
